$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = -0.3425203228057133
    "C2"  = 0.1370775318625127
    "D2"  = 1.367787377993229
    "B3"  = -0.0860618221422374
    "C3"  = 1.383680499826445
    "B4"  = 1.32438096272888
    "B5"  = 1.231661280012998
    "C5"  = 0.1923542655531081
    "D5"  = 0.2120950120634008
    "E5"  = 0.3792844388692188
    "B6"  = 0.4395685548510502
    "C6"  = 0.3295146744469067
    "D6"  = 0.2496258682164595
    "B7"  = 0.4940016924669799
    "C7"  = 0.2435019605816055
    "B8"  = 0.3242193037695071
    "B9"  = 0.394042449657095
    "C9"  = 0.3427404628639549
    "D9"  = 0.2242390996078211
    "E9"  = 0.055485660899395
    "B10" = 0.5010464375566571
    "C10" = 0.2455732575174918
    "D10" = 0.0184750902009912
    "B11" = 0.4070253497240054
    "C11" = 0.0588001744469144
    "B12" = 0.1842008206034934
    "B13" = -0.3665930774731743
    "C13" = 0.101671561665663
    "D13" = 0.04785045814007061
    "B14" = -0.1157322460162614
    "C14" = 0.09802926010601931
    "B15" = 0.2279386281717039
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
